{"js": "// Update the sheet date header and all 24 three-digit / one-digit division\n// prompts in the table to the next day's auto-generated values.\n//\n// Every \"old\" string below occurs exactly once in the document body (the\n// date line plus 24 distinct \"NNN\u00f7N=\" cells out of the worksheet's 40\n// table cells, the rest being intentionally blank answer rows), so a\n// literal, case-sensitive whole-document search-and-replace unambiguously\n// targets the right run without disturbing its formatting (font/size are\n// left untouched because only the run's text is rewritten in place).\nconst replacements = [\n  [\"2024-06-20 Thursday\", \"2024-06-21 Friday\"],\n  [\"497\u00f75=\", \"989\u00f76=\"],\n  [\"134\u00f78=\", \"786\u00f77=\"],\n  [\"889\u00f77=\", \"228\u00f73=\"],\n  [\"417\u00f79=\", \"291\u00f76=\"],\n  [\"384\u00f79=\", \"744\u00f74=\"],\n  [\"372\u00f77=\", \"922\u00f78=\"],\n  [\"684\u00f74=\", \"649\u00f78=\"],\n  [\"237\u00f79=\", \"710\u00f76=\"],\n  [\"986\u00f75=\", \"178\u00f79=\"],\n  [\"490\u00f75=\", \"977\u00f72=\"],\n  [\"882\u00f76=\", \"784\u00f78=\"],\n  [\"614\u00f72=\", \"784\u00f74=\"],\n  [\"216\u00f73=\", \"210\u00f77=\"],\n  [\"980\u00f74=\", \"531\u00f72=\"],\n  [\"971\u00f78=\", \"803\u00f76=\"],\n  [\"500\u00f72=\", \"491\u00f76=\"],\n  [\"342\u00f74=\", \"182\u00f78=\"],\n  [\"832\u00f77=\", \"946\u00f72=\"],\n  [\"255\u00f78=\", \"260\u00f72=\"],\n  [\"199\u00f79=\", \"679\u00f78=\"],\n  [\"553\u00f75=\", \"174\u00f72=\"],\n  [\"578\u00f73=\", \"769\u00f79=\"],\n  [\"948\u00f72=\", \"939\u00f78=\"],\n  [\"787\u00f76=\", \"113\u00f76=\"],\n  [\"325\u00f79=\", \"121\u00f78=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `expected exactly 1 match for ${JSON.stringify(oldText)}, found ${results.items.length}`\n    );\n  }\n\n  for (const foundRange of results.items) {\n    foundRange.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and all 24 three-digit / one-digit\n# division prompts to the next day's generated values. Each pair is\n# [oldText, newText]; old text strings are unique within the document,\n# so an exact Find/Replace (no wildcards, case-sensitive) is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-06-20 Thursday\", \"2024-06-21 Friday\"),\n  @(\"497\u00f75=\", \"989\u00f76=\"),\n  @(\"134\u00f78=\", \"786\u00f77=\"),\n  @(\"889\u00f77=\", \"228\u00f73=\"),\n  @(\"417\u00f79=\", \"291\u00f76=\"),\n  @(\"384\u00f79=\", \"744\u00f74=\"),\n  @(\"372\u00f77=\", \"922\u00f78=\"),\n  @(\"684\u00f74=\", \"649\u00f78=\"),\n  @(\"237\u00f79=\", \"710\u00f76=\"),\n  @(\"986\u00f75=\", \"178\u00f79=\"),\n  @(\"490\u00f75=\", \"977\u00f72=\"),\n  @(\"882\u00f76=\", \"784\u00f78=\"),\n  @(\"614\u00f72=\", \"784\u00f74=\"),\n  @(\"216\u00f73=\", \"210\u00f77=\"),\n  @(\"980\u00f74=\", \"531\u00f72=\"),\n  @(\"971\u00f78=\", \"803\u00f76=\"),\n  @(\"500\u00f72=\", \"491\u00f76=\"),\n  @(\"342\u00f74=\", \"182\u00f78=\"),\n  @(\"832\u00f77=\", \"946\u00f72=\"),\n  @(\"255\u00f78=\", \"260\u00f72=\"),\n  @(\"199\u00f79=\", \"679\u00f78=\"),\n  @(\"553\u00f75=\", \"174\u00f72=\"),\n  @(\"578\u00f73=\", \"769\u00f79=\"),\n  @(\"948\u00f72=\", \"939\u00f78=\"),\n  @(\"787\u00f76=\", \"113\u00f76=\"),\n  @(\"325\u00f79=\", \"121\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  # wdFindContinue=1, wdReplaceAll=2 -- MatchCase=$true, MatchWholeWord=$false,\n  # MatchWildcards=$false so the division sign / literal text matches exactly.\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"replace failed: '$oldText' was not found in the document\"\n  }\n}\n\n"}
